# Weekly price update for "Hortaliza, Comercializadora del Agro de Limarí - Tomate".
# A new week's worth of data (2 rows: "Larga vida" Segunda/Tercera) is inserted
# at the top of the price-history block (row 312), pushing all the existing
# history rows down by two rows (312-368 -> 314-370).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the historical block; Excel shifts every
# row from 312 down to 370 automatically (rows 312-368 -> 314-370) and carries
# the date-column's number format along with it.
$ws.Rows("312:313").Insert()

# --- New row 312: Tomate / Larga vida / Segunda ---
$ws.Range("A312").Value = 2
$ws.Range("B312").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C312").Value = "Coquimbo"
$ws.Range("D312").Value = 44476
$ws.Range("E312").Value = 4
$ws.Range("F312").Value = 100112020
$ws.Range("G312").Value = "Tomate"
$ws.Range("H312").Value = "Larga vida"
$ws.Range("I312").Value = "Segunda"
$ws.Range("J312").Value = 1800
$ws.Range("K312").Value = 12500
$ws.Range("L312").Value = 13000
$ws.Range("M312").Value = 12750
$ws.Range("N312").Value = "`$/bandeja 18 kilos"
$ws.Range("O312").Value = "Provincia de Limarí"
$ws.Range("P312").Value = 708
$ws.Range("Q312").Value = 18
$ws.Range("R312").Value = "Hortaliza"

# --- New row 313: Tomate / Larga vida / Tercera ---
$ws.Range("A313").Value = 2
$ws.Range("B313").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C313").Value = "Coquimbo"
$ws.Range("D313").Value = 44476
$ws.Range("E313").Value = 4
$ws.Range("F313").Value = 100112020
$ws.Range("G313").Value = "Tomate"
$ws.Range("H313").Value = "Larga vida"
$ws.Range("I313").Value = "Tercera"
$ws.Range("J313").Value = 1500
$ws.Range("K313").Value = 10500
$ws.Range("L313").Value = 11000
$ws.Range("M313").Value = 10750
$ws.Range("N313").Value = "`$/bandeja 18 kilos"
$ws.Range("O313").Value = "Provincia de Limarí"
$ws.Range("P313").Value = 597
$ws.Range("Q313").Value = 18
$ws.Range("R313").Value = "Hortaliza"
